$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-24 Wednesday", "2025-12-25 Thursday"),
    @("49×58=2842", "86×23=1978"),
    @("46×70=3220", "20×29=580"),
    @("77×70=5390", "96×59=5664"),
    @("65×16=1040", "89×22=1958"),
    @("72×62=4464", "85×50=4250"),
    @("61×51=3111", "87×90=7830"),
    @("95×99=9405", "17×30=510"),
    @("73×71=5183", "99×58=5742"),
    @("12×71=852", "15×80=1200"),
    @("55×22=1210", "98×47=4606"),
    @("24×98=2352", "93×99=9207"),
    @("21×42=882", "97×32=3104"),
    @("66×69=4554", "15×95=1425"),
    @("40×68=2720", "26×28=728"),
    @("43×42=1806", "53×58=3074"),
    @("50×19=950", "14×37=518"),
    @("74×62=4588", "59×86=5074"),
    @("55×85=4675", "26×70=1820"),
    @("93×66=6138", "96×22=2112"),
    @("30×81=2430", "64×98=6272"),
    @("27×69=1863", "70×47=3290"),
    @("29×96=2784", "99×41=4059"),
    @("22×14=308", "66×90=5940"),
    @("62×12=744", "32×59=1888"),
    @("41×39=1599", "70×79=5530")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
